# Increase MaxInvest Storage Adapt Szenarios Existing Units
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")

# MaxInvest (column S) increased from 8 to 15 for BESS units in rows 7-11
$ws.Range("S7:S11").Value = 15

# ExisUnits (column E) increased from 28 to 33 for BESS7 (row 10)
$ws.Range("E10").Value = 33

# Update the active cell selection to reflect where the user left off editing
$ws.Range("L22").Select()
